$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B3").Value = 1720.481741820667
$ws.Range("D7").Value = 1705.024607732537
